$wb = $excel.ActiveWorkbook

# --- Add the new "OrgData" worksheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$orgData = $wb.Worksheets.Add($null, $lastSheet)
$orgData.Name = "OrgData"

# Header row - reuses the existing "OrgName" shared string (same text as on the
# Organization sheet), plus an empty second row that will later be filled in
# dynamically with the created Organisation's name.
$orgData.Range("A1").Value = "OrgName"
$orgData.Range("A2").Value = ""

# Give the two populated cells a thin box border (new cell style referenced
# from the updated styles.xml).
$orgData.Range("A1:A2").Borders.LineStyle = 1

# Match the column sizing used on the other sheets.
$orgData.Columns.Item(1).ColumnWidth = 16

# The new sheet becomes the active / selected tab, with A2 as the active cell.
[void]$orgData.Range("A2").Select()

# --- Organization sheet: selection moves, it is no longer the selected tab ---
$org = $wb.Worksheets.Item("Organization")
[void]$org.Range("C15").Select()

# --- Opportunities sheet: selection moves back to A2 ---
$opp = $wb.Worksheets.Item("Opportunities")
[void]$opp.Range("A2").Select()

# Re-activate the OrgData sheet last so it stays the active sheet/tab in the
# saved workbook (activeTab / tabSelected), with A2 selected.
[void]$orgData.Activate()
[void]$orgData.Range("A2").Select()
